$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StateCounters")

# Update the MD SPA_Last counter (B2) from 9405 to 9412
$ws.Range("B2").Value = 9412
